$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '70.720.94'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -2.58%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.629.14'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '584.40'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -2.34%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '175.93'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.86%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.636'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +4.32%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '3.621.89'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('E9').Value = '  +0.00%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.196'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -5.81%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '6.80'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +15.98%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.617'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.26%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '48.45'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -4.04%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000283'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.74%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '4.216.19'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.46%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '674.21'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -4.68%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '9.01'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.11%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '3.634.04'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +1.66%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '70.791.14'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -2.58%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.123'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('E21').Value = '  -4.43%  '
$ws.Range('E22').Value = '  -2.10%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.945'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.66%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '17.19'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -4.65%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '100.01'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -4.53%  '
$ws.Range('E26').Value = '  -3.30%  '
$ws.Range('E27').Value = '  -3.08%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  -2.25%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '34.59'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -2.52%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '9.14'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.37%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.27'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -6.22%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '7.56'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.50%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.39'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -6.29%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.96'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -5.88%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '573.28'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -3.43%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '11.08'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -2.49%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.107'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.44%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '58.46'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.60%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.09%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0452'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.547.03'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.77%  '
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('E44').Value = '  -4.03%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '34.34'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -4.91%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0₃0732'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -6.62%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.68'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -4.78%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +3.70%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.136'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +2.28%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '137.93'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.39%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.88'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -5.03%  '
